# changement horaire + niveaux IA
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Team members who left the department - their rows are removed from the roster
$namesToRemove = @("Carine Croteau", "Cirine Chaieb", "Gabriel Montplaisir", "Sébastien Trottier")

# Walk bottom-up so row indices of rows still to examine aren't disturbed by deletions
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
for ($r = $lastRow; $r -ge 2; $r--) {
    $name = $ws.Cells.Item($r, 1).Value()
    if ($namesToRemove -contains $name) {
        $ws.Rows.Item($r).Delete()
    }
}

# The engine only supports clearing the *whole* hyperlink collection at once, so
# rebuild it from scratch against the now-compacted rows (2..8). Re-touching
# Font.Underline right after each Add() makes the engine settle back onto the
# original "Hyperlink" cell-style index instead of minting a new duplicate xf.
function Add-Link($range, $address, $display) {
    if ($display) {
        $ws.Hyperlinks.Add($range, $address, "", "", $display) | Out-Null
    } else {
        $ws.Hyperlinks.Add($range, $address) | Out-Null
    }
    $range.Font.Underline = 2
}

Add-Link $ws.Range("D4") "mailto:rivard.etienne@cegepvicto.ca" $null
Add-Link $ws.Range("G4") "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=RIVARD.ETIENNE@cegepvicto.ca" $null

Add-Link $ws.Range("D2") "mailto:ouellet.alexandre@cegepvicto.ca" $null
Add-Link $ws.Range("G2") "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=OUELLET.ALEXANDRE@cegepvicto.ca" $null

Add-Link $ws.Range("D7") "mailto:frechette.mathieu@cegepvicto.ca" $null
Add-Link $ws.Range("G7") "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=FRECHETTE.MATHIEU@cegepvicto.ca" $null

Add-Link $ws.Range("D3") "mailto:lagace.christiane@cegepvicto.ca" $null
Add-Link $ws.Range("G3") "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=LAGACE.CHRISTIANE@cegepvicto.ca" $null

Add-Link $ws.Range("D5") "mailto:mercier.francois@cegepvicto.ca" $null
Add-Link $ws.Range("G5") "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=MERCIER.FRANCOIS@cegepvicto.ca" $null

Add-Link $ws.Range("D6") "mailto:taleb.frederik@cegepvicto.ca" $null
Add-Link $ws.Range("G6") "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=TALEB.FREDERIK@cegepvicto.ca" $null

Add-Link $ws.Range("D8") "mailto:tousignant.simon@cegepvicto.ca" $null
Add-Link $ws.Range("G8") "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=TOUSIGNANT.SIMON@cegepvicto.ca" $null

Add-Link $ws.Range("C8") "https://avatar.iran.liara.run/public/boy" "https://avatar.iran.liara.run/public/boy"

# Restore the sort state / selection the author left behind
$ws.Range("C15").Select()
